$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp label (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Junio de 2020 a las 09:49"

# --- Row 4: Estados Unidos (no rank change, values updated) ---
$ws.Range("B4").Value = 2553068
$ws.Range("C4").Value = 112
$ws.Range("D4").Value = 1068768
$ws.Range("E4").Value = 1356660

# --- Row 7: India (no rank change, values updated) ---
$ws.Range("B7").Value = 509737
$ws.Range("C7").Value = 291
$ws.Range("D7").Value = 296028
$ws.Range("E7").Value = 198009
$ws.Range("G7").Value = 11
$ws.Range("H7").Value = 15700

# --- Row 36: Singapur (no rank change, values updated) ---
$ws.Range("B36").Value = 43246
$ws.Range("C36").Value = 291
$ws.Range("E36").Value = 6395

# --- Rows 37/38: Ucrania overtakes Irak in ranking ---
# Row 37 becomes Ucrania with fresh data
$ws.Range("A37").Value = "Ucrania"
$ws.Range("B37").Value = 42065
$ws.Range("C37").Value = 948
$ws.Range("D37").Value = 18701
$ws.Range("E37").Value = 22254
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 24
$ws.Range("H37").Value = 1110

# Row 38 becomes Irak with its previous (unchanged) data
$ws.Range("A38").Value = "Irak"
$ws.Range("B38").Value = 41193
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 18859
$ws.Range("E38").Value = 20775
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 1559

# --- Row 44: Afganistan (no rank change, values updated) ---
$ws.Range("B44").Value = 30616
$ws.Range("C44").Value = 165
$ws.Range("D44").Value = 10674
$ws.Range("E44").Value = 19239
$ws.Range("G44").Value = 20
$ws.Range("H44").Value = 703

# --- Rows 51/52: Armenia overtakes Nigeria in ranking ---
# Row 51 becomes Armenia with fresh data
$ws.Range("A51").Value = "Armenia"
$ws.Range("B51").Value = 23909
$ws.Range("C51").Value = 662
$ws.Range("D51").Value = 12911
$ws.Range("E51").Value = 10577
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 11
$ws.Range("H51").Value = 421

# Row 52 becomes Nigeria with its previous (unchanged) data
$ws.Range("A52").Value = "Nigeria"
$ws.Range("B52").Value = 23298
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 8253
$ws.Range("E52").Value = 14491
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 554

# --- Row 92: Hungria (no rank change, values updated) ---
$ws.Range("B92").Value = 4138
$ws.Range("C92").Value = 11
$ws.Range("D92").Value = 2681
$ws.Range("E92").Value = 879

# --- Row 110: Estonia (no rank change, values updated) ---
$ws.Range("D110").Value = 1812
$ws.Range("E110").Value = 105

# --- Row 126: Letonia (no rank change, values updated) ---
$ws.Range("B126").Value = 1115
$ws.Range("C126").Value = 3
$ws.Range("E126").Value = 153

# --- Row 157: Surinam (no rank change, values updated) ---
$ws.Range("B157").Value = 391
$ws.Range("C157").Value = 2
$ws.Range("E157").Value = 200

# --- Rows 165/166: Guyana overtakes Mongolia in ranking ---
# Row 165 becomes Guyana with fresh data
$ws.Range("A165").Value = "Guyana"
$ws.Range("B165").Value = 230
$ws.Range("C165").Value = 15
$ws.Range("D165").Value = 109
$ws.Range("E165").Value = 109
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 12

# Row 166 becomes Mongolia with its previous (unchanged) data
$ws.Range("A166").Value = "Mongolia"
$ws.Range("B166").Value = 219
$ws.Range("C166").Value = 0
$ws.Range("D166").Value = 175
$ws.Range("E166").Value = 44
$ws.Range("F166").Value = 0
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 0

# --- Rows 201-204: cosmetic re-sort among tied countries (same totals) ---
$ws.Range("A201").Value = "Santa Lucia"
$ws.Range("A202").Value = "Laos"
$ws.Range("A203").Value = "Fiyi"
$ws.Range("A204").Value = "Dominica"

# --- Rows 208-209: cosmetic re-sort among tied countries (same totals) ---
$ws.Range("A208").Value = "Groenlandia"
$ws.Range("A209").Value = "Islas Malvinas"

# --- Rows 212/213: Seychelles overtakes Montserrat in ranking ---
# Row 212 becomes Seychelles with its previous (unchanged) data
$ws.Range("A212").Value = "Seychelles"
$ws.Range("D212").Value = 11
$ws.Range("H212").Value = 0

# Row 213 becomes Montserrat with its previous (unchanged) data
$ws.Range("A213").Value = "Montserrat"
$ws.Range("D213").Value = 10
$ws.Range("H213").Value = 1
